# Applies the LOM3107.xlsx content restructuring described in the commit diff:
#  - "Objetivos:", "Programa resumido:", "Programa:", "Metodo:", "Criterio:",
#    "Norma de recuperacao:" and "Bibliografia:" get their real descriptive text
#    (previously B/C of those rows held misplaced faculty names/labels).
#  - A new "Docentes responsaveis:" row is added, followed by the four faculty
#    members each on their own row (instead of one crammed per section row).
#  - The used range grows from A1:C22 to A1:C27 and several rows lose their
#    explicit 60/120 custom row height, falling back to the sheet default (15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Objetivos:" gets its real text in B10/C10 (A10 label/style/height stay put) ---
$ws.Range("B10").Value = 'Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões, deformações e deslocamentos em estruturas compostas por barras em regime elástico-linear sob carregamento axial, torção e flexão.Desenvolver aplicações práticas para dimensionamento de barras em condições de carregamentos mistos.Prover o conhecimento dos fenômenos de flambagem, com aplicações práticas para dimensionamento de colunas.Descrever a metodologia para análise dos estados planos de tensão e deformação, bem como a aplicação da lei de Hooke para casos multiaxiais.Apresentar conceitos básicos sobre energia de deformação.'
$ws.Range("C10").Value = 'Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões, deformações e deslocamentos em estruturas compostas por barras em regime elástico-linear sob carregamento axial, torção e flexão.Desenvolver aplicações práticas para dimensionamento de barras em condições de carregamentos mistos.Prover o conhecimento dos fenômenos de flambagem, com aplicações práticas para dimensionamento de colunas.Descrever a metodologia para análise dos estados planos de tensão e deformação, bem como a aplicação da lei de Hooke para casos multiaxiais.Apresentar conceitos básicos sobre energia de deformação.'

# --- Drop the old, misaligned rows 11-22 ---
$ws.Range("A11:A22").EntireRow.Delete()

# --- Insert 17 fresh rows (new 11-27); they inherit col A/B/C styles from row 10 above ---
$ws.Range("A11:A27").EntireRow.Insert()

# Row 11
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Clear()
$ws.Range("C11").Clear()
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Rows.Item(12).AutoFit()

# Row 13
$ws.Range("B13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("C13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("A13").Clear()
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("B14").Value = '3480026 - João Paulo Pascon'
$ws.Range("C14").Value = '3480026 - João Paulo Pascon'
$ws.Range("A14").Clear()
$ws.Rows.Item(14).AutoFit()

# Row 15
$ws.Range("B15").Value = '5840793 - Sérgio Schneider'
$ws.Range("C15").Value = '5840793 - Sérgio Schneider'
$ws.Range("A15").Clear()
$ws.Rows.Item(15).AutoFit()

# Row 16
$ws.Range("B16").Value = '7797767 - Viktor Pastoukhov'
$ws.Range("C16").Value = '7797767 - Viktor Pastoukhov'
$ws.Range("A16").Clear()
$ws.Rows.Item(16).AutoFit()

# Row 17
$ws.Range("A17").Value = 'Programa resumido:'
$ws.Range("B17").Value = 'Considerações fundamentais. Conceito de tensão. Conceito de deformação. Lei de Hooke. Carga Axial. Torção em barras de seção circular. Flexão em vigas isostáticas de seção simétrica. Cargas combinadas. Flambagem de colunas. Análise de Tensão e Deformação. Lei de Hooke Multiaxial. Energia de deformação.'
$ws.Range("C17").Value = 'Considerações fundamentais. Conceito de tensão. Conceito de deformação. Lei de Hooke. Carga Axial. Torção em barras de seção circular. Flexão em vigas isostáticas de seção simétrica. Cargas combinadas. Flambagem de colunas. Análise de Tensão e Deformação. Lei de Hooke Multiaxial. Energia de deformação.'
$ws.Rows.Item(17).RowHeight = 60

# Row 18
$ws.Range("A18").Value = 'Short syllabus:'
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = 'Programa:'
$ws.Range("B19").Value = '1. Considerações fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes.2. Conceito de tensão: Tensão Normal; Tensão Cisalhante; Tensões admissíveis.3. Conceito de deformação: Deformação Normal; Deformação por Cisalhamento.4. Lei de Hooke: Elasticidade linear e o Módulo de Young; Lei de Hooke para Cisalhamento.5. Carga Axial: Deslocamentos em sistemas isostáticos; Efeitos da Temperatura; Sistemas Hiperestáticos.6. Torção em barras de seção circular: Momento de inércia polar; Análise das tensões em eixos de seção maciça e seção vazada; Cálculo das rotações relativas entre seções adjacentes; Eixos estaticamente indeterminados; Torção e tração combinadas.7. Flexão em vigas isostáticas de seção simétrica: Forças concentradas e forças distribuídas; Diagramas de força cortante e momento fletor para uma viga carregada; Momento de inércia, eixos principais de inércia; Flexão em Vigas de Seção Simétrica; Determinação das Tensões Normais; Deflexões em vigas: equação diferencial da linha elástica; Tensões de cisalhamento em vigas. Tensões de cisalhamento em barras de paredes finas.8. Cargas combinadas: Modos Mistos de Carregamento. Projeto de barras submetidas a cargas axiais, transversais e torcionais.9. Flambagem de colunas: Raio de giração. Fórmula de Euler para colunas biarticuladas. Fatores de correção para outras condições de contorno. Projeto de colunas de aço e de outras ligas submetidas a um carregamento centrado.10. Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Transformação do Estado Plano de Deformação.11. Lei de Hooke Multiaxial: Elasticidade, Homogeneidade e Isotropia; Coeficiente de Poisson; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas.12. Energia de deformação: Densidade de energia de deformação. Energia de deformação elástica para tensões normais. Energia de deformação elástica para tensões de cisalhamento. Projeto para carregamento por impacto. Métodos de energia: teorema de Castigliano e suas aplicações.'
$ws.Range("C19").Value = '1. Considerações fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes.2. Conceito de tensão: Tensão Normal; Tensão Cisalhante; Tensões admissíveis.3. Conceito de deformação: Deformação Normal; Deformação por Cisalhamento.4. Lei de Hooke: Elasticidade linear e o Módulo de Young; Lei de Hooke para Cisalhamento.5. Carga Axial: Deslocamentos em sistemas isostáticos; Efeitos da Temperatura; Sistemas Hiperestáticos.6. Torção em barras de seção circular: Momento de inércia polar; Análise das tensões em eixos de seção maciça e seção vazada; Cálculo das rotações relativas entre seções adjacentes; Eixos estaticamente indeterminados; Torção e tração combinadas.7. Flexão em vigas isostáticas de seção simétrica: Forças concentradas e forças distribuídas; Diagramas de força cortante e momento fletor para uma viga carregada; Momento de inércia, eixos principais de inércia; Flexão em Vigas de Seção Simétrica; Determinação das Tensões Normais; Deflexões em vigas: equação diferencial da linha elástica; Tensões de cisalhamento em vigas. Tensões de cisalhamento em barras de paredes finas.8. Cargas combinadas: Modos Mistos de Carregamento. Projeto de barras submetidas a cargas axiais, transversais e torcionais.9. Flambagem de colunas: Raio de giração. Fórmula de Euler para colunas biarticuladas. Fatores de correção para outras condições de contorno. Projeto de colunas de aço e de outras ligas submetidas a um carregamento centrado.10. Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Transformação do Estado Plano de Deformação.11. Lei de Hooke Multiaxial: Elasticidade, Homogeneidade e Isotropia; Coeficiente de Poisson; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas.12. Energia de deformação: Densidade de energia de deformação. Energia de deformação elástica para tensões normais. Energia de deformação elástica para tensões de cisalhamento. Projeto para carregamento por impacto. Métodos de energia: teorema de Castigliano e suas aplicações.'
$ws.Rows.Item(19).RowHeight = 120

# Row 20
$ws.Range("A20").Value = 'Syllabus:'
$ws.Range("B20").Clear()
$ws.Range("C20").Clear()
$ws.Rows.Item(20).RowHeight = 120

# Row 21
$ws.Range("A21").Value = 'Avaliação:'
$ws.Range("B21").Clear()
$ws.Range("C21").Clear()
$ws.Rows.Item(21).AutoFit()

# Row 22
$ws.Range("A22").Value = 'Método:'
$ws.Range("B22").Value = 'Os alunos serão avaliados por meio de três conjuntos de notas: duas provas escritas (P1 e P2) envolvendo o conteúdo teórico ministrado em sala de aula; exercícios (EX) propostos regularmente para serem entregues e discutidos na aula subsequente; e seminários (SM) em grupo ao final da disciplina.'
$ws.Range("C22").Value = 'Os alunos serão avaliados por meio de três conjuntos de notas: duas provas escritas (P1 e P2) envolvendo o conteúdo teórico ministrado em sala de aula; exercícios (EX) propostos regularmente para serem entregues e discutidos na aula subsequente; e seminários (SM) em grupo ao final da disciplina.'
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = 'Critério:'
$ws.Range("B23").Value = 'Nota Final (NF) = 70%((P1+P2)/2)+20%(EX)+10%(SM).'
$ws.Range("C23").Value = 'Nota Final (NF) = 70%((P1+P2)/2)+20%(EX)+10%(SM).'
$ws.Rows.Item(23).RowHeight = 60

# Row 24
$ws.Range("A24").Value = 'Norma de recuperação:'
$ws.Range("B24").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C24").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Rows.Item(24).RowHeight = 60

# Row 25
$ws.Range("A25").Value = 'Bibliografia:'
$ws.Range("B25").Value = '1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p.3. R.R. CRAIG, Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p.4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p.5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p.6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p.7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais. Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p.'
$ws.Range("C25").Value = '1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p.3. R.R. CRAIG, Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p.4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p.5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p.6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p.7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais. Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p.'
$ws.Rows.Item(25).RowHeight = 120

# Row 26
$ws.Range("A26").Value = 'Requisitos:'
$ws.Range("B26").Clear()
$ws.Range("C26").Clear()
$ws.Rows.Item(26).AutoFit()

# Row 27
$ws.Range("B27").Value = "LOM3099 -  Estática  (Requisito fraco)`n"
$ws.Range("C27").Value = "LOM3099 -  Estática  (Requisito fraco)`n"
$ws.Range("A27").Clear()
$ws.Rows.Item(27).RowHeight = 30
